# Updates the crypto price/volume table (and two ranking swaps) on Sheet1
# to reflect a refreshed data pull, per the source diff.
#
# Cells whose new value could be misread as a number by Excel (e.g. "1.00",
# "0.0000201", "9.50") are written through Set-TextCellValue, which temporarily
# forces a text number format so the literal string is preserved, then restores
# the cell to its original (unstyled/General) appearance.

function Set-TextCellValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).NumberFormat = "General"
    $ws.Range($cellRef).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "91.721.11"
$ws.Range("E2").Value = "  +1.59%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.126.17"
$ws.Range("E3").Value = "  +2.26%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.58%  "

# Row 5: Solana
Set-TextCellValue $ws "D5" "246.23"

# Row 6: BNB
Set-TextCellValue $ws "D6" "618.44"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7: XRP
$ws.Range("E7").Value = "  -1.41%  "

# Row 8: Dogecoin
$ws.Range("E8").Value = "  +5.88%  "

# Row 9: USDC
Set-TextCellValue $ws "D9" "1.00"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10: LidoStakedEther
$ws.Range("D10").Value = "3.121.61"
$ws.Range("E10").Value = "  +2.14%  "

# Row 11: Cardano
Set-TextCellValue $ws "D11" "0.737"
$ws.Range("E11").Value = "  +0.47%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +2.05%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  +2.20%  "

# Row 14: Avalanche
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCellValue $ws "D14" "34.87"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15: Toncoin
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCellValue $ws "D15" "5.59"
$ws.Range("E15").Value = "  +3.05%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "91.573.45"
$ws.Range("E16").Value = "  +1.40%  "

# Row 17: WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.704.86"
$ws.Range("E17").Value = "  +1.67%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.089.43"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19: SuiNetwork
$ws.Range("E19").Value = "  +1.24%  "

# Row 20: Chainlink
Set-TextCellValue $ws "D20" "14.97"
$ws.Range("E20").Value = "  +4.57%  "

# Row 21: Polkadot
Set-TextCellValue $ws "D21" "5.87"
$ws.Range("E21").Value = "  +2.34%  "

# Row 22: Uniswap
Set-TextCellValue $ws "D22" "9.50"
$ws.Range("E22").Value = "  +5.96%  "

# Row 23: BitcoinCash
Set-TextCellValue $ws "D23" "449.27"
$ws.Range("E23").Value = "  +2.46%  "

# Row 24: PEPE
Set-TextCellValue $ws "D24" "0.0000201"
$ws.Range("E24").Value = "  -3.52%  "

# Row 25: NEARProtocol
Set-TextCellValue $ws "D25" "5.90"
$ws.Range("E25").Value = "  +6.03%  "

# Row 26: Litecoin
$ws.Range("E26").Value = "  +5.68%  "

# Row 27: Aptos
Set-TextCellValue $ws "D27" "11.79"
$ws.Range("E27").Value = "  +0.58%  "

# Row 29: Hedera
Set-TextCellValue $ws "D29" "0.147"
$ws.Range("E29").Value = "  +32.25%  "

# Row 30: Dai
$ws.Range("E30").Value = "  +0.21%  "

# Row 31: Stellar
$ws.Range("E31").Value = "  -3.26%  "

# Row 32: Cronos
$ws.Range("E32").Value = "  -7.48%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCellValue $ws "D33" "9.39"
$ws.Range("E33").Value = "  +3.45%  "

# Row 34: Kaspa
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCellValue $ws "D34" "0.176"
$ws.Range("E34").Value = "  +6.75%  "

# Row 35: Binance-PegBSC-USD
$ws.Range("E35").Value = "  -0.33%  "

# Row 36: RenderToken
Set-TextCellValue $ws "D36" "7.92"
$ws.Range("E36").Value = "  +3.79%  "

# Row 37: EthereumClassic
Set-TextCellValue $ws "D37" "26.30"
$ws.Range("E37").Value = "  +0.23%  "

# Row 38: MantraDAO
$ws.Range("E38").Value = "  +0.44%  "

# Row 39: PancakeSwap
$ws.Range("E39").Value = "  +2.52%  "

# Row 40: Bittensor
Set-TextCellValue $ws "D40" "492.45"
$ws.Range("E40").Value = "  +0.69%  "

# Row 41: Fetch.AI
$ws.Range("E41").Value = "  +2.00%  "

# Row 42: PolygonEcosystemToken
Set-TextCellValue $ws "D42" "0.445"
$ws.Range("E42").Value = "  +7.35%  "

# Row 43: dogwifhat
Set-TextCellValue $ws "D43" "3.40"
$ws.Range("E43").Value = "  -4.55%  "

# Row 44: WhiteBITCoin
$ws.Range("E44").Value = "  +0.30%  "

# Row 45: USDe
$ws.Range("E45").Value = "  -0.04%  "

# Row 46: Monero
Set-TextCellValue $ws "D46" "158.70"
$ws.Range("E46").Value = "  +3.10%  "

# Row 47: ARBITRUM
Set-TextCellValue $ws "D47" "0.709"
$ws.Range("E47").Value = "  +4.62%  "

# Row 48: Stacks
$ws.Range("E48").Value = "  +2.37%  "

# Row 49: ImmutableX
$ws.Range("E49").Value = "  +3.43%  "

# Row 50: Filecoin
Set-TextCellValue $ws "D50" "4.43"
$ws.Range("E50").Value = "  +0.85%  "

# Row 51: OKB
Set-TextCellValue $ws "D51" "44.05"
$ws.Range("E51").Value = "  -0.05%  "
